$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6324274688181788
$ws.Range("C2").Value = 0.1329156646247043
$ws.Range("E2").Value = 0.5094825150310243
$ws.Range("F2").Value = 1.996243267601514
$ws.Range("G2").Value = 0.00240438185066641
$ws.Range("J2").Value = 0.03188770868450375
$ws.Range("M2").Value = 0.5145993241915932
$ws.Range("N2").Value = 1.107572114958224
$ws.Range("O2").Value = 1.99336100911799

$ws.Range("B3").Value = 0.558104782886403
$ws.Range("C3").Value = 0.1179959648263491
$ws.Range("E3").Value = 0.505442239682921
$ws.Range("F3").Value = 1.981317332101412
$ws.Range("G3").Value = 0.00240709800672441
$ws.Range("J3").Value = 0.0323068925796961
$ws.Range("M3").Value = 0.4850356947614358
$ws.Range("N3").Value = 1.118262505974016
$ws.Range("O3").Value = 1.991779750431164

$ws.Range("B4").Value = 0.5124463963037726
$ws.Range("C4").Value = 0.1087840884221123
$ws.Range("E4").Value = 0.503175650801353
$ws.Range("F4").Value = 1.973336956450709
$ws.Range("G4").Value = 0.002408854717538593
$ws.Range("J4").Value = 0.03257854721530729
$ws.Range("M4").Value = 0.4670614605206964
$ws.Range("N4").Value = 1.125300794649533
$ws.Range("O4").Value = 1.992338727010747

$ws.Range("B5").Value = 0.4938350888542402
$ws.Range("C5").Value = 0.1050175068498618
$ws.Range("E5").Value = 0.5023059051203731
$ws.Range("F5").Value = 1.97038247495783
$ws.Range("G5").Value = 0.002409593033793302
$ws.Range("J5").Value = 0.03269283806876055
$ws.Range("M5").Value = 0.4597818810343313
$ws.Range("N5").Value = 1.128288272875729
$ws.Range("O5").Value = 1.992950430492755

$ws.Range("B6").Value = 0.4907444107247443
$ws.Range("C6").Value = 0.1043913093225228
$ws.Range("E6").Value = 0.502164742080879
$ws.Range("F6").Value = 1.969909853491671
$ws.Range("G6").Value = 0.002409716988040499
$ws.Range("J6").Value = 0.03271203277069068
$ws.Range("M6").Value = 0.4585758437679317
$ws.Range("N6").Value = 1.128791549177592
$ws.Range("O6").Value = 1.993075168059107

$ws.Range("B7").Value = 0.5121954173627046
$ws.Range("C7").Value = 0.1087333420253742
$ws.Range("E7").Value = 0.5031637027621372
$ws.Range("F7").Value = 1.97329590662217
$ws.Range("G7").Value = 0.002408864583776125
$ws.Range("J7").Value = 0.03258007404811636
$ws.Range("M7").Value = 0.4669631026682168
$ws.Range("N7").Value = 1.125340601655189
$ws.Range("O7").Value = 1.992345423241119

$ws.Range("B8").Value = 0.6068064833755216
$ws.Range("C8").Value = 0.1277820701216967
$ws.Range("E8").Value = 0.5080450111583374
$ws.Range("F8").Value = 1.990850895752487
$ws.Range("G8").Value = 0.002405299954310378
$ws.Range("J8").Value = 0.03202928079197731
$ws.Range("M8").Value = 0.504369031579941
$ws.Range("N8").Value = 1.111159734900347
$ws.Range("O8").Value = 1.992497786721231

$ws.Range("B9").Value = 0.7921188745432914
$ws.Range("C9").Value = 0.1647244371420982
$ws.Range("E9").Value = 0.5193150524101
$ws.Range("F9").Value = 2.034685755430417
$ws.Range("G9").Value = 0.002399012649979289
$ws.Range("J9").Value = 0.03106241755664962
$ws.Range("M9").Value = 0.5791240061617984
$ws.Range("N9").Value = 1.087113123636541
$ws.Range("O9").Value = 2.004974159292658

$ws.Range("B10").Value = 0.9281077717318453
$ws.Range("C10").Value = 0.1916085370599205
$ws.Range("E10").Value = 0.5286298197521404
$ws.Range("F10").Value = 2.072652942185471
$ws.Range("G10").Value = 0.002394817551390237
$ws.Range("J10").Value = 0.03042105294008746
$ws.Range("M10").Value = 0.6348939340501261
$ws.Range("N10").Value = 1.071737358208679
$ws.Range("O10").Value = 2.021622754679782

$ws.Range("B11").Value = 0.9899331654424373
$ws.Range("C11").Value = 0.2037817517349083
$ws.Range("E11").Value = 0.5330920529230951
$ws.Range("F11").Value = 2.091182220308283
$ws.Range("G11").Value = 0.002393000276698403
$ws.Range("J11").Value = 0.03014425060216119
$ws.Range("M11").Value = 0.6604479765349538
$ws.Range("N11").Value = 1.065239704572491
$ws.Range("O11").Value = 2.030834203096447

$ws.Range("B12").Value = 1.013338863821843
$ws.Range("C12").Value = 0.2083831526624067
$ws.Range("E12").Value = 0.5348140968254427
$ws.Range("F12").Value = 2.098379994382839
$ws.Range("G12").Value = 0.002392325152008766
$ws.Range("J12").Value = 0.03004158377994415
$ws.Range("M12").Value = 0.6701508570850763
$ws.Range("N12").Value = 1.062850651938952
$ws.Range("O12").Value = 2.034558821774937

$ws.Range("B13").Value = 1.008298320642552
$ws.Range("C13").Value = 0.2073925324679351
$ws.Range("E13").Value = 0.5344417883617183
$ws.Range("F13").Value = 2.09682176450805
$ws.Range("G13").Value = 0.002392469973179874
$ws.Range("J13").Value = 0.03006359917520651
$ws.Range("M13").Value = 0.6680600119410371
$ws.Range("N13").Value = 1.06336199824959
$ws.Range("O13").Value = 2.033746128108191

$ws.Range("B14").Value = 0.9918588968561153
$ws.Range("C14").Value = 0.204160479913071
$ws.Range("E14").Value = 0.533233079578487
$ws.Range("F14").Value = 2.091770753696267
$ws.Range("G14").Value = 0.002392944472894702
$ws.Range("J14").Value = 0.03013576099937598
$ws.Range("M14").Value = 0.6612457168495922
$ws.Range("N14").Value = 1.065041723321777
$ws.Range("O14").Value = 2.031135885821641

$ws.Range("B15").Value = 0.9817884342620573
$ws.Range("C15").Value = 0.2021796632930375
$ws.Range("E15").Value = 0.53249691466727
$ws.Range("F15").Value = 2.088700461246233
$ws.Range("G15").Value = 0.002393236813695188
$ws.Range("J15").Value = 0.03018024252227303
$ws.Range("M15").Value = 0.6570751556836854
$ws.Range("N15").Value = 1.066079911321239
$ws.Range("O15").Value = 2.029567855987011

$ws.Range("B16").Value = 0.9240665167969837
$ws.Range("C16").Value = 0.1908118342648493
$ws.Range("E16").Value = 0.5283427238162162
$ws.Range("F16").Value = 2.07146734347171
$ws.Range("G16").Value = 0.002394938142366064
$ws.Range("J16").Value = 0.03043944372214025
$ws.Range("M16").Value = 0.6332275909405212
$ws.Range("N16").Value = 1.07217199483383
$ws.Range("O16").Value = 2.021053801979519

$ws.Range("B17").Value = 0.8886459166025134
$ws.Range("C17").Value = 0.1838234214803265
$ws.Range("E17").Value = 0.525851830532801
$ws.Range("F17").Value = 2.061217729543401
$ws.Range("G17").Value = 0.002396005140278953
$ws.Range("J17").Value = 0.03060228720106828
$ws.Range("M17").Value = 0.6186447706485438
$ws.Range("N17").Value = 1.07603656102178
$ws.Range("O17").Value = 2.016250865273946

$ws.Range("B18").Value = 0.8682695333733932
$ws.Range("C18").Value = 0.1797985658215282
$ws.Range("E18").Value = 0.5244403059147018
$ws.Range("F18").Value = 2.055440798142342
$ws.Range("G18").Value = 0.002396627427447174
$ws.Range("J18").Value = 0.03069735826455755
$ws.Range("M18").Value = 0.6102744596329188
$ws.Range("N18").Value = 1.078306122499995
$ws.Range("O18").Value = 2.013642457265433

$ws.Range("B19").Value = 0.8613698872900954
$ws.Range("C19").Value = 0.1784349139569485
$ws.Range("E19").Value = 0.5239660255730527
$ws.Range("F19").Value = 2.053505151434692
$ws.Range("G19").Value = 0.002396839598301691
$ws.Range("J19").Value = 0.03072978944506932
$ws.Range("M19").Value = 0.6074434050333366
$ws.Range("N19").Value = 1.079082588464189
$ws.Range("O19").Value = 2.012785736587688

$ws.Range("B20").Value = 0.892416857920864
$ws.Range("C20").Value = 0.1845679005491832
$ws.Range("E20").Value = 0.5261147995664643
$ws.Range("F20").Value = 2.062296565100269
$ws.Range("G20").Value = 0.00239589066936671
$ws.Range("J20").Value = 0.0305848065230494
$ws.Range("M20").Value = 0.62019534446965
$ws.Range("N20").Value = 1.075620331196504
$ws.Range("O20").Value = 2.016746189232293

$ws.Range("B21").Value = 0.9966877310744735
$ws.Range("C21").Value = 0.2051100398494441
$ws.Range("E21").Value = 0.5335872305961686
$ws.Range("F21").Value = 2.093249439831723
$ws.Range("G21").Value = 0.002392804747550548
$ws.Range("J21").Value = 0.03011450690103068
$ws.Range("M21").Value = 0.6632465343581657
$ws.Range("N21").Value = 1.06454640755225
$ws.Range("O21").Value = 2.031896152988622

$ws.Range("B22").Value = 1.064797976508714
$ws.Range("C22").Value = 0.2184868588421409
$ws.Range("E22").Value = 0.5386591033862658
$ws.Range("F22").Value = 2.114534871595083
$ws.Range("G22").Value = 0.002390863884671657
$ws.Range("J22").Value = 0.02981968432505422
$ws.Range("M22").Value = 0.6915350602229182
$ws.Range("N22").Value = 1.057725539920135
$ws.Range("O22").Value = 2.043176124305518

$ws.Range("B23").Value = 1.028449942776604
$ws.Range("C23").Value = 0.2113519190590978
$ws.Range("E23").Value = 0.5359349438696341
$ws.Range("F23").Value = 2.103077729272229
$ws.Range("G23").Value = 0.002391892829318791
$ws.Range("J23").Value = 0.02997588832548637
$ws.Range("M23").Value = 0.6764231388899589
$ws.Range("N23").Value = 1.061327840231726
$ws.Range("O23").Value = 2.03702934856716

$ws.Range("B24").Value = 0.8907120545516705
$ws.Range("C24").Value = 0.1842313437380767
$ws.Range("E24").Value = 0.5259958473346629
$ws.Range("F24").Value = 2.061808463137538
$ws.Range("G24").Value = 0.002395942394037638
$ws.Range("J24").Value = 0.03059270502006228
$ws.Range("M24").Value = 0.6194942878236844
$ws.Range("N24").Value = 1.07580835989436
$ws.Range("O24").Value = 2.016521777169459

$ws.Range("B25").Value = 0.7420131980952647
$ws.Range("C25").Value = 0.1547753148178685
$ws.Range("E25").Value = 0.5160844851631836
$ws.Range("F25").Value = 2.021817285799628
$ws.Range("G25").Value = 0.002400638729743013
$ws.Range("J25").Value = 0.03131186131075392
$ws.Range("M25").Value = 0.5587514616831584
$ws.Range("N25").Value = 1.093215860162104
$ws.Range("O25").Value = 2.000289236793719
